$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I19").Value = 0.6842399319549618
$ws.Range("J19").Value = 0.2163398458498075
$ws.Range("K19").Value = 0.1645595860744188
$ws.Range("L19").Value = 2.635073052411303

$ws.Range("I20").Value = 0.9585908516801056
$ws.Range("J20").Value = 0.4914461501244363
$ws.Range("K20").Value = 0.4577586258096178
$ws.Range("L20").Value = 2.225721949624816
